$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.509.05"
$ws.Range("E2").Value = "'  -0.22%  "
$ws.Range("D3").Value = "'1.919.16"
$ws.Range("E3").Value = "'  -0.20%  "
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("D5").Value = "'245.82"
$ws.Range("E5").Value = "'  -0.62%  "
$ws.Range("E6").Value = "'  +0.06%  "
$ws.Range("D7").Value = "'0.4854"
$ws.Range("E7").Value = "'  +2.77%  "
$ws.Range("D8").Value = "'0.2895"
$ws.Range("E8").Value = "'  +0.18%  "
$ws.Range("D9").Value = "'0.06704"
$ws.Range("E9").Value = "'  -1.27%  "
$ws.Range("D10").Value = "'111.99"
$ws.Range("E10").Value = "'  +5.94%  "
$ws.Range("D11").Value = "'19.27"
$ws.Range("E11").Value = "'  +4.38%  "
$ws.Range("D12").Value = "'1.919.36"
$ws.Range("E12").Value = "'  +0.21%  "
$ws.Range("D13").Value = "'0.07590"
$ws.Range("E13").Value = "'  -1.23%  "
$ws.Range("D14").Value = "'5.358"
$ws.Range("E14").Value = "'  +1.11%  "
$ws.Range("D15").Value = "'0.6733"
$ws.Range("E15").Value = "'  +0.03%  "
$ws.Range("D16").Value = "'295.89"
$ws.Range("E16").Value = "'  +1.47%  "
$ws.Range("D17").Value = "'30.518.72"
$ws.Range("E17").Value = "'  -0.25%  "
$ws.Range("D18").Value = "'13.04"
$ws.Range("E18").Value = "'  +0.93%  "
$ws.Range("E19").Value = "'  +0.10%  "
$ws.Range("D20").Value = "'0.000007563"
$ws.Range("E20").Value = "'  -0.66%  "
$ws.Range("B21").Value = "'Uniswap"
$ws.Range("C21").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'5.526"
$ws.Range("E21").Value = "'  +0.46%  "
$ws.Range("B22").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "'2.169.94"
$ws.Range("E22").Value = "'  +0.23%  "
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "'  +0.13%  "
$ws.Range("D24").Value = "'6.451"
$ws.Range("E24").Value = "'  +1.88%  "
$ws.Range("E25").Value = "'  +0.64%  "
$ws.Range("D26").Value = "'164.14"
$ws.Range("E26").Value = "'  -2.82%  "
$ws.Range("D27").Value = "'20.23"
$ws.Range("E27").Value = "'  -4.80%  "
$ws.Range("D28").Value = "'2.097"
$ws.Range("E28").Value = "'  -1.16%  "
$ws.Range("D29").Value = "'0.1076"
$ws.Range("E29").Value = "'  -0.11%  "
$ws.Range("D30").Value = "'1.438"
$ws.Range("E30").Value = "'  +3.07%  "
$ws.Range("D31").Value = "'4.131"
$ws.Range("E31").Value = "'  -1.32%  "
$ws.Range("D32").Value = "'4.084"
$ws.Range("E32").Value = "'  -1.25%  "
$ws.Range("D33").Value = "'0.05019"
$ws.Range("E33").Value = "'  -1.01%  "
$ws.Range("D34").Value = "'0.7390"
$ws.Range("E34").Value = "'  -0.64%  "
$ws.Range("D35").Value = "'1.139"
$ws.Range("E35").Value = "'  -1.27%  "
$ws.Range("D36").Value = "'1.000"
$ws.Range("D37").Value = "'2.719"
$ws.Range("E37").Value = "'  -1.10%  "
$ws.Range("D38").Value = "'0.02020"
$ws.Range("E38").Value = "'  -2.71%  "
$ws.Range("D39").Value = "'2.697"
$ws.Range("E39").Value = "'  +0.07%  "
$ws.Range("D40").Value = "'109.99"
$ws.Range("E40").Value = "'  -1.37%  "
$ws.Range("D41").Value = "'2.010"
$ws.Range("E41").Value = "'  -2.42%  "
$ws.Range("D42").Value = "'0.4434"
$ws.Range("E42").Value = "'  -0.14%  "
$ws.Range("D43").Value = "'0.8643"
$ws.Range("E43").Value = "'  -2.03%  "
$ws.Range("E44").Value = "'  -0.62%  "
$ws.Range("D45").Value = "'70.15"
$ws.Range("E45").Value = "'  +4.53%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "'  +0.05%  "
$ws.Range("D47").Value = "'7.230"
$ws.Range("E47").Value = "'  -0.40%  "
$ws.Range("D48").Value = "'48.81"
$ws.Range("E48").Value = "'  -0.40%  "
$ws.Range("D49").Value = "'9.187"
$ws.Range("E49").Value = "'  -1.97%  "
$ws.Range("D50").Value = "'0.1224"
$ws.Range("E50").Value = "'  -1.19%  "
$ws.Range("D51").Value = "'0.2507"
$ws.Range("E51").Value = "'  +2.05%  "

# Re-normalize style: the quote-prefix text entry above marks the cell
# xf with quotePrefix=1; reset back to the workbook default "Normal" style
# so the written cells keep their original (unstyled) formatting.
$ws.Range("B2:E51").Style = "Normal"

